# Rename the worksheet tabs to match their corresponding class names,
# then reorder the trailing "beamline" sheets so Radiation (-> laserTab)
# moves to the very end, after LBeamline (-> laserTransportTab), and make
# the former EBeamline (-> bunchTransportTab) the active sheet.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("UndC").Name       = "undulator"
$wb.Worksheets.Item("BeamP").Name      = "particleBeam"
$wb.Worksheets.Item("Precision").Name  = "precisionSRW"
$wb.Worksheets.Item("Wavefront").Name  = "wavefrontSRW"
$wb.Worksheets.Item("Radiation").Name  = "laserTab"
$wb.Worksheets.Item("EBeamline").Name  = "bunchTransportTab"
$wb.Worksheets.Item("LBeamline").Name  = "laserTransportTab"

# Move "laserTab" (formerly Radiation) to the end of the tab strip, i.e.
# right after "laserTransportTab" (formerly LBeamline). This leaves the
# tab order as:
#   undulator, particleBeam, precisionSRW, wavefrontSRW,
#   bunchTransportTab, laserTransportTab, laserTab
$wsLaserTab = $wb.Worksheets.Item("laserTab")
$wsLaserTransportTab = $wb.Worksheets.Item("laserTransportTab")
$wsLaserTab.Move($null, $wsLaserTransportTab)

# The active/selected tab becomes "bunchTransportTab" (formerly EBeamline).
$wb.Worksheets.Item("bunchTransportTab").Activate()
